$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The updated TPM-based analysis collapses the previous 4 data rows into 3
# (the FAPs/MuSCs <-> MuSCs row is dropped) and refreshes all of the
# numeric statistics. Remove the now-unused last row first.
$ws.Rows("5:5").Delete() | Out-Null

# Row 2: ECs -> Efna5/Epha5 -> MuSCs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efna5"
$ws.Cells.Item(2,3).Value = "Epha5"
$ws.Cells.Item(2,4).Value = "MuSCs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.114918
$ws.Cells.Item(2,8).Value = 0.344754
$ws.Cells.Item(2,9).Value = 0.04640425382421802
$ws.Cells.Item(2,10).Value = 0.04640425382421801
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.01136166666666667
$ws.Cells.Item(2,14).Value = 0.034085
$ws.Cells.Item(2,15).Value = 1
$ws.Cells.Item(2,16).Value = 1
$ws.Cells.Item(2,17).Value = 0.00130566001
$ws.Cells.Item(2,18).Value = 0.01175094009
$ws.Cells.Item(2,19).Value = 0.04640425382421802
$ws.Cells.Item(2,20).Value = 0.04640425382421801

# Row 3: FAPs -> Efna5/Epha5 -> MuSCs
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Efna5"
$ws.Cells.Item(3,3).Value = "Epha5"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.030023666666667
$ws.Cells.Item(3,8).Value = 6.090071
$ws.Cells.Item(3,9).Value = 0.819730011809897
$ws.Cells.Item(3,10).Value = 0.819730011809897
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.01136166666666667
$ws.Cells.Item(3,14).Value = 0.034085
$ws.Cells.Item(3,15).Value = 1
$ws.Cells.Item(3,16).Value = 1
$ws.Cells.Item(3,17).Value = 0.02306445222611111
$ws.Cells.Item(3,18).Value = 0.207580070035
$ws.Cells.Item(3,19).Value = 0.819730011809897
$ws.Cells.Item(3,20).Value = 0.819730011809897

# Row 4: MuSCs -> Efna5/Epha5 -> MuSCs
$ws.Cells.Item(4,1).Value = "MuSCs"
$ws.Cells.Item(4,2).Value = "Efna5"
$ws.Cells.Item(4,3).Value = "Epha5"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.3315123333333334
$ws.Cells.Item(4,8).Value = 0.994537
$ws.Cells.Item(4,9).Value = 0.133865734365885
$ws.Cells.Item(4,10).Value = 0.133865734365885
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.01136166666666667
$ws.Cells.Item(4,14).Value = 0.034085
$ws.Cells.Item(4,15).Value = 1
$ws.Cells.Item(4,16).Value = 1
$ws.Cells.Item(4,17).Value = 0.003766532627222222
$ws.Cells.Item(4,18).Value = 0.03389879364499999
$ws.Cells.Item(4,19).Value = 0.133865734365885
$ws.Cells.Item(4,20).Value = 0.133865734365885
